$d = $word.ActiveDocument

# 1) Replace the opening "Sou " with the new lead-in sentence about
#    10 years of freelance PC/network experience.
$d.Content.Find.Execute(
    "Sou ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Com 10 anos de experiência em montagem e manutenção de computadores e configurações de redes domésticas e empresariais como profissional autônomo, também possuo ",
    2
)

# 2) Drop the now-redundant "um profissional com " before "formação em
#    Educação Física" in the following run.
$d.Content.Find.Execute(
    "um profissional com formação em Educação Física",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "formação em Educação Física",
    2
)

# 3) Remove the blank paragraph that used to separate the summary
#    paragraph from the "Pós-Graduações" heading.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq [char]13) {
        $n = $p.Next()
        if ($n -ne $null -and $n.Range.Text.StartsWith("Pós-Gradua")) {
            $r = $d.Range($p.Range.Start, $n.Range.Start)
            $r.Delete()
            break
        }
    }
}

# 4) Merge the four runs of "Formações acadêmicas" heading into a single
#    run of text. Delete the trailing three runs' text ("es acadêmicas")
#    and re-insert it at the end of the first run ("Formaçõ") so that
#    run's (empty) rPr is preserved rather than dropped by a wholesale
#    text replacement.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Formaçõ" + "es")) {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $tail = $d.Range($pStart + 7, $pEnd - 1)
        $tail.Delete()
        $insertionPoint = $d.Range($pStart + 7, $pStart + 7)
        $insertionPoint.InsertAfter("es acadêmicas")
        break
    }
}

# 5) Merge "bacharelado" + " em Administração..." into one run.
$d.Content.Find.Execute(
    "bacharelado em Administração pela Faculdade Anhanguera (2023 - 2026) – Em Andamento",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "bacharelado em Administração pela Faculdade Anhanguera (2023 - 2026) – Em Andamento",
    2
)
